$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated figures from the 2020-08-10 Fonds de solidarite data refresh.
# Column C = nombre_aides, Column D = montant_total (kept as text to preserve exact decimal formatting).
$updates = @(
    @{ Row = 20; C = "165"; D = "420599.00" }
    @{ Row = 21; C = "323"; D = "1117929.00" }
    @{ Row = 22; C = "156"; D = "441137.26" }
    @{ Row = 28; C = "245"; D = "627174.00" }
    @{ Row = 30; C = "497"; D = "1945299.70" }
    @{ Row = 32; C = "356"; D = "1153722.37" }
    @{ Row = 33; C = "10"; D = "31000.00" }
    @{ Row = 40; C = "109"; D = "277077.69" }
    @{ Row = 41; C = "72"; D = "338959.98" }
    @{ Row = 42; C = "112"; D = "461688.99" }
    @{ Row = 45; C = "319"; D = "867006.74" }
    @{ Row = 47; C = "562"; D = "2122430.99" }
    @{ Row = 48; C = "375"; D = "1230567.16" }
    @{ Row = 51; C = "3280"; D = "7467348.36" }
    @{ Row = 53; C = "3781"; D = "12753819.79" }
    @{ Row = 55; C = "3866"; D = "11788438.49" }
    @{ Row = 73; C = "364"; D = "893635.70" }
    @{ Row = 75; C = "878"; D = "2935419.89" }
    @{ Row = 76; C = "497"; D = "1602502.87" }
    @{ Row = 91; C = "540"; D = "1321218.67" }
    @{ Row = 92; C = "8"; D = "22200.00" }
    @{ Row = 93; C = "1052"; D = "3502629.98" }
    @{ Row = 95; C = "952"; D = "2819531.31" }
    @{ Row = 96; C = "14"; D = "43500.00" }
    @{ Row = 97; C = "45"; D = "166833.00" }
)

foreach ($u in $updates) {
    $cRange = $ws.Range("C" + $u.Row)
    $dRange = $ws.Range("D" + $u.Row)
    $cRange.NumberFormat = "@"
    $dRange.NumberFormat = "@"
    $cRange.Value = $u.C
    $dRange.Value = $u.D
}
